$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its original plain-text representation.
# Some updated price strings (e.g. "324.34") would otherwise be auto-detected
# by Excel as numbers, which would corrupt values like "28.923.77" style
# thousand-separated text and introduce floating point artifacts.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.923.77"
$ws.Range("E2").Value = "  -1.94%  "
$ws.Range("D3").Value = "1.900.75"
$ws.Range("E3").Value = "  -4.01%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "324.34"
$ws.Range("E5").Value = "  -1.04%  "
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("D7").Value = "0.4579"
$ws.Range("E7").Value = "  -1.68%  "
$ws.Range("D8").Value = "0.3808"
$ws.Range("E8").Value = "  -2.69%  "
$ws.Range("E9").Value = "  -3.12%  "
$ws.Range("D10").Value = "0.9737"
$ws.Range("E10").Value = "  -1.93%  "
$ws.Range("D11").Value = "21.97"
$ws.Range("E11").Value = "  -3.98%  "
$ws.Range("D12").Value = "1.902.14"
$ws.Range("E12").Value = "  -4.96%  "
$ws.Range("D13").Value = "6.913"
$ws.Range("E13").Value = "  -3.87%  "
$ws.Range("D14").Value = "5.641"
$ws.Range("E14").Value = "  -3.34%  "
$ws.Range("D15").Value = "0.07033"
$ws.Range("E15").Value = "  -0.84%  "
$ws.Range("D16").Value = "1.005"
$ws.Range("E16").Value = "  -0.22%  "
$ws.Range("D17").Value = "83.56"
$ws.Range("E17").Value = "  -4.70%  "
$ws.Range("D18").Value = "0.000009461"
$ws.Range("E18").Value = "  -5.10%  "
$ws.Range("D19").Value = "16.58"
$ws.Range("E19").Value = "  -4.08%  "
$ws.Range("E20").Value = "  -0.12%  "
$ws.Range("D21").Value = "28.897.89"
$ws.Range("E21").Value = "  -2.08%  "
$ws.Range("D22").Value = "5.275"
$ws.Range("E22").Value = "  -4.97%  "
$ws.Range("D23").Value = "10.84"
$ws.Range("E23").Value = "  -3.13%  "
$ws.Range("D24").Value = "2.094"
$ws.Range("E24").Value = "  -0.72%  "
$ws.Range("D25").Value = "158.09"
$ws.Range("E25").Value = "  -0.22%  "
$ws.Range("D26").Value = "19.00"
$ws.Range("E26").Value = "  -3.16%  "
$ws.Range("D27").Value = "5.593"
$ws.Range("E27").Value = "  -3.74%  "
$ws.Range("D28").Value = "117.09"
$ws.Range("E28").Value = "  -2.08%  "
$ws.Range("D29").Value = "1.830"
$ws.Range("E29").Value = "  -3.94%  "
$ws.Range("D30").Value = "0.09235"
$ws.Range("E30").Value = "  -2.03%  "
$ws.Range("D31").Value = "0.8557"
$ws.Range("E31").Value = "  -3.99%  "
$ws.Range("D32").Value = "5.068"
$ws.Range("E32").Value = "  -3.17%  "
$ws.Range("E33").Value = "  -6.80%  "
$ws.Range("D34").Value = "2.997"
$ws.Range("E34").Value = "  -6.23%  "
$ws.Range("D35").Value = "0.05654"
$ws.Range("E35").Value = "  -2.74%  "
$ws.Range("D36").Value = "1.139"
$ws.Range("E36").Value = "  -3.22%  "
$ws.Range("D37").Value = "1.003"
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("D38").Value = "0.02039"
$ws.Range("E38").Value = "  -2.78%  "
$ws.Range("D39").Value = "0.5468"
$ws.Range("E39").Value = "  -4.60%  "
$ws.Range("D40").Value = "7.374"
$ws.Range("E40").Value = "  -5.32%  "
$ws.Range("D41").Value = "0.1749"
$ws.Range("E41").Value = "  -3.03%  "
$ws.Range("D42").Value = "9.261"
$ws.Range("E42").Value = "  -4.38%  "
$ws.Range("D43").Value = "2.754"
$ws.Range("E43").Value = "  -1.51%  "
$ws.Range("D44").Value = "0.5140"
$ws.Range("E44").Value = "  -4.21%  "
$ws.Range("D45").Value = "11.20"
$ws.Range("E45").Value = "  -5.39%  "
$ws.Range("D46").Value = "0.06808"
$ws.Range("E46").Value = "  -1.82%  "
$ws.Range("B47").Value = "PEPE"
$ws.Range("C47").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D47").Value = "0.000002607"
$ws.Range("E47").Value = "  -15.15%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "2.058"
$ws.Range("E48").Value = "  -5.52%  "
$ws.Range("D49").Value = "109.82"
$ws.Range("E49").Value = "  -3.86%  "
$ws.Range("E50").Value = "  -3.52%  "
$ws.Range("D51").Value = "1.002"
